$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 52
$templateRow = 51

# Copy the formatting (styles) of the last existing row onto the new row
# before touching any values, so the new row inherits the same per-column
# cell styles (s="4" for most columns, s="5" for the hyperlink columns).
$ws.Range("A$templateRow`:J$templateRow").Copy()
$ws.Range("A$newRow`:J$newRow").PasteSpecial(-4122)

# Fill in the new source record (Auto Bild).
$ws.Cells.Item($newRow, 1).Value  = "Q_AB"
$ws.Cells.Item($newRow, 2).Value  = "Auto Bild"
$ws.Cells.Item($newRow, 3).Value  = "Auto Bild"
$ws.Cells.Item($newRow, 4).Value  = "Auto Bild"
$ws.Cells.Item($newRow, 5).Value  = "Auto Bild"
$ws.Cells.Item($newRow, 6).Value  = "www.autobild.de"
$ws.Cells.Item($newRow, 7).Value  = "www.autobild.de"
$ws.Cells.Item($newRow, 10).Value = "ab"

# Add hyperlinks on the Homepage De / Homepage En cells.
$ws.Hyperlinks.Add($ws.Range("F$newRow"), "https://www.autobild.de", "", "", "www.autobild.de") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G$newRow"), "https://www.autobild.de", "", "", "www.autobild.de") | Out-Null

# Adding a hyperlink re-styles the cell with the built-in "Hyperlink" style,
# so re-apply the source formatting for those two cells afterwards.
$ws.Range("F$templateRow`:G$templateRow").Copy()
$ws.Range("F$newRow`:G$newRow").PasteSpecial(-4122)
